# Re-apply a built-in PowerPoint table style to the table on slide 5
# (was previously using the custom/imported Google Slides table style
# "{683C7126-05E7-4A18-B424-4520985EC973}"; switch it to the built-in
# style "{BACBF737-9DFD-4612-B60F-2B9C808B9C37}").

$p = $ppt.ActivePresentation

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            $table.ApplyStyle("{BACBF737-9DFD-4612-B60F-2B9C808B9C37}", $true)
        }
    }
}
